# Add season record columns (Wins, Losses, Ties) to the OAK_2012 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new headers in AD1:AF1, matching the existing header style.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows (2-52): every player row gets the team's season record.
for ($row = 2; $row -le 52; $row++) {
    $ws.Cells.Item($row, 30).Value = 94  # AD
    $ws.Cells.Item($row, 31).Value = 68  # AE
    $ws.Cells.Item($row, 32).Value = 0   # AF
}
